$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, shifting existing rows 132..261 down to 133..262
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new data record
$ws.Range("A132").Value = 10
$ws.Range("B132").Value = "Vega Modelo de Temuco"
$ws.Range("C132").Value = "La Araucanía"
$ws.Range("D132").Value = 44586
$ws.Range("E132").Value = 9
$ws.Range("F132").Value = 100112009
$ws.Range("G132").Value = "Acelga"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 65
$ws.Range("K132").Value = 7000
$ws.Range("L132").Value = 8000
$ws.Range("M132").Value = 7462
$ws.Range("N132").Value = "$/docena de atados (12 kilos)"
$ws.Range("O132").Value = "Provincia de Cautín"
$ws.Range("P132").Value = 622
$ws.Range("Q132").Value = 12
$ws.Range("R132").Value = "Hortaliza"

# Match the style (date format) used on column D for the rest of the rows
$ws.Range("D132").NumberFormat = $ws.Range("D133").NumberFormat
